# Update quarter date format: strip the leading "Qtr" from the period
# labels in column A (e.g. "Qtr1-2010" -> "1-2010", "Qtr2-2014" -> "2-2014").
#
# The rows are touched in the same order Excel's own Find & Replace (Ctrl+H,
# Replace All on "Qtr" -> "") walked them, so the rebuilt shared-string
# table comes out in the same order as the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowOrder = @(2,3,6,7,4,5,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42)

foreach ($r in $rowOrder) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "Qtr*") {
        $cell.Value = $val -replace "^Qtr", ""
    }
}

# Restore the view state recorded the last time the sheet was saved:
# scrolled down so row 13 is at the top, with A44 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("A44").Select()
